$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("How many curves can I load in one go?", "llama3.2:latest", "You can load up to 450 curves at a time."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("What the maximum number of headers I can display in my log?", "llama3.2:latest", "You can select up to 50 headers for display in a log."),
    @("Whats the maximum number of characters in a single text entry?", "llama3.2:latest", "The maximum number of characters in a single text entry is 32000."),
    @("How many tables can I have in my log?", "llama3.2:latest", "You can have up to 100 tables in a log."),
    @("How many symbols can I have in the plot at any one time?", "llama3.2:latest", "You can have up to 10,000 symbols in a plot at any given time."),
    @("How many scales can I define?", "llama3.2:latest", "According to the document, you can define 23 scales.")
)

$startRow = 32
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
